# Update math-fact strings in the division-practice table.
$d = $word.ActiveDocument

$replacements = @(
    @("33÷3=", "15÷7="),
    @("12÷6=", "23÷7="),
    @("53÷6=", "21÷5="),
    @("66÷6=", "18÷3="),
    @("36÷7=", "28÷8="),
    @("31÷7=", "26÷2="),
    @("67÷7=", "89÷6="),
    @("47÷4=", "88÷8="),
    @("49÷3=", "20÷5="),
    @("74÷3=", "24÷7="),
    @("25÷2=", "56÷9="),
    @("59÷4=", "28÷6="),
    @("35÷2=", "59÷4="),
    @("68÷7=", "80÷8="),
    @("59÷3=", "44÷5="),
    @("53÷8=", "26÷4="),
    @("81÷8=", "98÷6="),
    @("43÷7=", "62÷6="),
    @("90÷7=", "18÷4="),
    @("45÷6=", "14÷3="),
    @("61÷2=", "35÷7="),
    @("77÷3=", "82÷6="),
    @("57÷8=", "56÷8="),
    @("79÷9=", "62÷5="),
    @("83÷4=", "57÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
